$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A20 value with the new, more precise serial timestamp
$ws.Cells.Item(20, 1).Value = 45876.75020091435

# Append new row 21 with updated readings
$ws.Cells.Item(21, 1).Value = 45876.79187491632
$ws.Cells.Item(21, 2).Value = 2025
$ws.Cells.Item(21, 3).Value = 28
$ws.Cells.Item(21, 4).Value = 14.83
$ws.Cells.Item(21, 5).Value = 88.45
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 6.24
$ws.Cells.Item(21, 8).Value = "ESE"
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = "19:00:17"

# Match the style of column A used for other date rows
$ws.Cells.Item(21, 1).NumberFormat = $ws.Cells.Item(20, 1).NumberFormat
